# Regenerate the "K" column (column G) values for the lyles_jordan save_data sheet.
# The commit replaces the previous Strike# derived figures with the new K values,
# recalculated as part of "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 9
    3  = 0
    4  = 6
    5  = 5
    6  = 1
    7  = 2
    8  = 2
    9  = 7
    10 = 9
    11 = 1
    12 = 7
    13 = 4
    14 = 2
    15 = 5
    16 = 5
    17 = 4
    18 = 7
    19 = 4
    20 = 4
    21 = 7
    22 = 1
    23 = 3
    24 = 3
    25 = 8
    26 = 8
    27 = 6
    28 = 5
    29 = 6
    30 = 3
    31 = 6
    32 = 4
    33 = 2
    34 = 2
    35 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
